$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 1).Value = "ECs"
$ws.Cells.Item(2, 2).Value = "Gnai2"
$ws.Cells.Item(2, 3).Value = "P2ry12"
$ws.Cells.Item(2, 4).Value = "ECs"
$ws.Cells.Item(2, 5).Value = 3
$ws.Cells.Item(2, 6).Value = 1
$ws.Cells.Item(2, 7).Value = 201.4397426666667
$ws.Cells.Item(2, 8).Value = 604.3192280000001
$ws.Cells.Item(2, 9).Value = 0.4833500233086392
$ws.Cells.Item(2, 10).Value = 0.4833500233086393
$ws.Cells.Item(2, 11).Value = 2
$ws.Cells.Item(2, 12).Value = 0.6666666666666666
$ws.Cells.Item(2, 13).Value = 0.04910833333333333
$ws.Cells.Item(2, 14).Value = 0.147325
$ws.Cells.Item(2, 15).Value = 0.002188553694087003
$ws.Cells.Item(2, 16).Value = 0.002188553694087002
$ws.Cells.Item(2, 17).Value = 9.892370029455556
$ws.Cells.Item(2, 18).Value = 89.0313302651
$ws.Cells.Item(2, 19).Value = 0.001057837479049161
$ws.Cells.Item(2, 20).Value = 0.001057837479049161

$ws.Cells.Item(3, 1).Value = "ECs"
$ws.Cells.Item(3, 2).Value = "Gnai2"
$ws.Cells.Item(3, 3).Value = "P2ry12"
$ws.Cells.Item(3, 4).Value = "MuSCs"
$ws.Cells.Item(3, 5).Value = 3
$ws.Cells.Item(3, 6).Value = 1
$ws.Cells.Item(3, 7).Value = 201.4397426666667
$ws.Cells.Item(3, 8).Value = 604.3192280000001
$ws.Cells.Item(3, 9).Value = 0.4833500233086392
$ws.Cells.Item(3, 10).Value = 0.4833500233086393
$ws.Cells.Item(3, 11).Value = 1
$ws.Cells.Item(3, 12).Value = 0.3333333333333333
$ws.Cells.Item(3, 13).Value = 0.2380986666666667
$ws.Cells.Item(3, 14).Value = 0.714296
$ws.Cells.Item(3, 15).Value = 0.01061106498877699
$ws.Cells.Item(3, 16).Value = 0.01061106498877699
$ws.Cells.Item(3, 17).Value = 47.96253414260978
$ws.Cells.Item(3, 18).Value = 431.6628072834881
$ws.Cells.Item(3, 19).Value = 0.005128858509654844
$ws.Cells.Item(3, 20).Value = 0.005128858509654844

$ws.Cells.Item(4, 1).Value = "ECs"
$ws.Cells.Item(4, 2).Value = "Gnai2"
$ws.Cells.Item(4, 3).Value = "P2ry12"
$ws.Cells.Item(4, 4).Value = "Resolving-Mac"
$ws.Cells.Item(4, 5).Value = 3
$ws.Cells.Item(4, 6).Value = 1
$ws.Cells.Item(4, 7).Value = 201.4397426666667
$ws.Cells.Item(4, 8).Value = 604.3192280000001
$ws.Cells.Item(4, 9).Value = 0.4833500233086392
$ws.Cells.Item(4, 10).Value = 0.4833500233086393
$ws.Cells.Item(4, 11).Value = 3
$ws.Cells.Item(4, 12).Value = 1
$ws.Cells.Item(4, 13).Value = 22.15150833333333
$ws.Cells.Item(4, 14).Value = 66.45452499999999
$ws.Cells.Item(4, 15).Value = 0.987200381317136
$ws.Cells.Item(4, 16).Value = 0.987200381317136
$ws.Cells.Item(4, 17).Value = 4462.194138345189
$ws.Cells.Item(4, 18).Value = 40159.7472451067
$ws.Cells.Item(4, 19).Value = 0.4771633273199353
$ws.Cells.Item(4, 20).Value = 0.4771633273199353

$ws.Cells.Item(5, 1).Value = "FAPs"
$ws.Cells.Item(5, 2).Value = "Gnai2"
$ws.Cells.Item(5, 3).Value = "P2ry12"
$ws.Cells.Item(5, 4).Value = "ECs"
$ws.Cells.Item(5, 5).Value = 3
$ws.Cells.Item(5, 6).Value = 1
$ws.Cells.Item(5, 7).Value = 65.41736466666667
$ws.Cells.Item(5, 8).Value = 196.252094
$ws.Cells.Item(5, 9).Value = 0.1569674599353791
$ws.Cells.Item(5, 10).Value = 0.1569674599353792
$ws.Cells.Item(5, 11).Value = 2
$ws.Cells.Item(5, 12).Value = 0.6666666666666666
$ws.Cells.Item(5, 13).Value = 0.04910833333333333
$ws.Cells.Item(5, 14).Value = 0.147325
$ws.Cells.Item(5, 15).Value = 0.002188553694087003
$ws.Cells.Item(5, 16).Value = 0.002188553694087002
$ws.Cells.Item(5, 17).Value = 3.212537749838889
$ws.Cells.Item(5, 18).Value = 28.91283974855
$ws.Cells.Item(5, 19).Value = 0.0003435317142930277
$ws.Cells.Item(5, 20).Value = 0.0003435317142930276

$ws.Cells.Item(6, 1).Value = "FAPs"
$ws.Cells.Item(6, 2).Value = "Gnai2"
$ws.Cells.Item(6, 3).Value = "P2ry12"
$ws.Cells.Item(6, 4).Value = "MuSCs"
$ws.Cells.Item(6, 5).Value = 3
$ws.Cells.Item(6, 6).Value = 1
$ws.Cells.Item(6, 7).Value = 65.41736466666667
$ws.Cells.Item(6, 8).Value = 196.252094
$ws.Cells.Item(6, 9).Value = 0.1569674599353791
$ws.Cells.Item(6, 10).Value = 0.1569674599353792
$ws.Cells.Item(6, 11).Value = 1
$ws.Cells.Item(6, 12).Value = 0.3333333333333333
$ws.Cells.Item(6, 13).Value = 0.2380986666666667
$ws.Cells.Item(6, 14).Value = 0.714296
$ws.Cells.Item(6, 15).Value = 0.01061106498877699
$ws.Cells.Item(6, 16).Value = 0.01061106498877699
$ws.Cells.Item(6, 17).Value = 15.57578730398045
$ws.Cells.Item(6, 18).Value = 140.182085735824
$ws.Cells.Item(6, 19).Value = 0.001665591918497557
$ws.Cells.Item(6, 20).Value = 0.001665591918497557

$ws.Cells.Item(7, 1).Value = "FAPs"
$ws.Cells.Item(7, 2).Value = "Gnai2"
$ws.Cells.Item(7, 3).Value = "P2ry12"
$ws.Cells.Item(7, 4).Value = "Resolving-Mac"
$ws.Cells.Item(7, 5).Value = 3
$ws.Cells.Item(7, 6).Value = 1
$ws.Cells.Item(7, 7).Value = 65.41736466666667
$ws.Cells.Item(7, 8).Value = 196.252094
$ws.Cells.Item(7, 9).Value = 0.1569674599353791
$ws.Cells.Item(7, 10).Value = 0.1569674599353792
$ws.Cells.Item(7, 11).Value = 3
$ws.Cells.Item(7, 12).Value = 1
$ws.Cells.Item(7, 13).Value = 22.15150833333333
$ws.Cells.Item(7, 14).Value = 66.45452499999999
$ws.Cells.Item(7, 15).Value = 0.987200381317136
$ws.Cells.Item(7, 16).Value = 0.987200381317136
$ws.Cells.Item(7, 17).Value = 1449.093298558372
$ws.Cells.Item(7, 18).Value = 13041.83968702535
$ws.Cells.Item(7, 19).Value = 0.1549583363025886
$ws.Cells.Item(7, 20).Value = 0.1549583363025886

$ws.Cells.Item(8, 1).Value = "MuSCs"
$ws.Cells.Item(8, 2).Value = "Gnai2"
$ws.Cells.Item(8, 3).Value = "P2ry12"
$ws.Cells.Item(8, 4).Value = "ECs"
$ws.Cells.Item(8, 5).Value = 3
$ws.Cells.Item(8, 6).Value = 1
$ws.Cells.Item(8, 7).Value = 60.43484133333334
$ws.Cells.Item(8, 8).Value = 181.304524
$ws.Cells.Item(8, 9).Value = 0.1450120099461104
$ws.Cells.Item(8, 10).Value = 0.1450120099461104
$ws.Cells.Item(8, 11).Value = 2
$ws.Cells.Item(8, 12).Value = 0.6666666666666666
$ws.Cells.Item(8, 13).Value = 0.04910833333333333
$ws.Cells.Item(8, 14).Value = 0.147325
$ws.Cells.Item(8, 15).Value = 0.002188553694087003
$ws.Cells.Item(8, 16).Value = 0.002188553694087002
$ws.Cells.Item(8, 17).Value = 2.967854333144444
$ws.Cells.Item(8, 18).Value = 26.7106889983
$ws.Cells.Item(8, 19).Value = 0.000317366570054541
$ws.Cells.Item(8, 20).Value = 0.000317366570054541

$ws.Cells.Item(9, 1).Value = "MuSCs"
$ws.Cells.Item(9, 2).Value = "Gnai2"
$ws.Cells.Item(9, 3).Value = "P2ry12"
$ws.Cells.Item(9, 4).Value = "MuSCs"
$ws.Cells.Item(9, 5).Value = 3
$ws.Cells.Item(9, 6).Value = 1
$ws.Cells.Item(9, 7).Value = 60.43484133333334
$ws.Cells.Item(9, 8).Value = 181.304524
$ws.Cells.Item(9, 9).Value = 0.1450120099461104
$ws.Cells.Item(9, 10).Value = 0.1450120099461104
$ws.Cells.Item(9, 11).Value = 1
$ws.Cells.Item(9, 12).Value = 0.3333333333333333
$ws.Cells.Item(9, 13).Value = 0.2380986666666667
$ws.Cells.Item(9, 14).Value = 0.714296
$ws.Cells.Item(9, 15).Value = 0.01061106498877699
$ws.Cells.Item(9, 16).Value = 0.01061106498877699
$ws.Cells.Item(9, 17).Value = 14.38945514167822
$ws.Cells.Item(9, 18).Value = 129.505096275104
$ws.Cells.Item(9, 19).Value = 0.001538731861691352
$ws.Cells.Item(9, 20).Value = 0.001538731861691352

$ws.Cells.Item(10, 1).Value = "MuSCs"
$ws.Cells.Item(10, 2).Value = "Gnai2"
$ws.Cells.Item(10, 3).Value = "P2ry12"
$ws.Cells.Item(10, 4).Value = "Resolving-Mac"
$ws.Cells.Item(10, 5).Value = 3
$ws.Cells.Item(10, 6).Value = 1
$ws.Cells.Item(10, 7).Value = 60.43484133333334
$ws.Cells.Item(10, 8).Value = 181.304524
$ws.Cells.Item(10, 9).Value = 0.1450120099461104
$ws.Cells.Item(10, 10).Value = 0.1450120099461104
$ws.Cells.Item(10, 11).Value = 3
$ws.Cells.Item(10, 12).Value = 1
$ws.Cells.Item(10, 13).Value = 22.15150833333333
$ws.Cells.Item(10, 14).Value = 66.45452499999999
$ws.Cells.Item(10, 15).Value = 0.987200381317136
$ws.Cells.Item(10, 16).Value = 0.987200381317136
$ws.Cells.Item(10, 17).Value = 1338.722891419011
$ws.Cells.Item(10, 18).Value = 12048.5060227711
$ws.Cells.Item(10, 19).Value = 0.1431559115143645
$ws.Cells.Item(10, 20).Value = 0.1431559115143645

$ws.Cells.Item(11, 1).Value = "Resolving-Mac"
$ws.Cells.Item(11, 2).Value = "Gnai2"
$ws.Cells.Item(11, 3).Value = "P2ry12"
$ws.Cells.Item(11, 4).Value = "ECs"
$ws.Cells.Item(11, 5).Value = 3
$ws.Cells.Item(11, 6).Value = 1
$ws.Cells.Item(11, 7).Value = 89.46554166666668
$ws.Cells.Item(11, 8).Value = 268.396625
$ws.Cells.Item(11, 9).Value = 0.2146705068098712
$ws.Cells.Item(11, 10).Value = 0.2146705068098712
$ws.Cells.Item(11, 11).Value = 2
$ws.Cells.Item(11, 12).Value = 0.6666666666666666
$ws.Cells.Item(11, 13).Value = 0.04910833333333333
$ws.Cells.Item(11, 14).Value = 0.147325
$ws.Cells.Item(11, 15).Value = 0.002188553694087003
$ws.Cells.Item(11, 16).Value = 0.002188553694087002
$ws.Cells.Item(11, 17).Value = 4.393503642013889
$ws.Cells.Item(11, 18).Value = 39.541532778125
$ws.Cells.Item(11, 19).Value = 0.0004698179306902728
$ws.Cells.Item(11, 20).Value = 0.0004698179306902727

$ws.Cells.Item(12, 1).Value = "Resolving-Mac"
$ws.Cells.Item(12, 2).Value = "Gnai2"
$ws.Cells.Item(12, 3).Value = "P2ry12"
$ws.Cells.Item(12, 4).Value = "MuSCs"
$ws.Cells.Item(12, 5).Value = 3
$ws.Cells.Item(12, 6).Value = 1
$ws.Cells.Item(12, 7).Value = 89.46554166666668
$ws.Cells.Item(12, 8).Value = 268.396625
$ws.Cells.Item(12, 9).Value = 0.2146705068098712
$ws.Cells.Item(12, 10).Value = 0.2146705068098712
$ws.Cells.Item(12, 11).Value = 1
$ws.Cells.Item(12, 12).Value = 0.3333333333333333
$ws.Cells.Item(12, 13).Value = 0.2380986666666667
$ws.Cells.Item(12, 14).Value = 0.714296
$ws.Cells.Item(12, 15).Value = 0.01061106498877699
$ws.Cells.Item(12, 16).Value = 0.01061106498877699
$ws.Cells.Item(12, 17).Value = 21.30162618344445
$ws.Cells.Item(12, 18).Value = 191.714635651
$ws.Cells.Item(12, 19).Value = 0.002277882698933237
$ws.Cells.Item(12, 20).Value = 0.002277882698933237

$ws.Cells.Item(13, 1).Value = "Resolving-Mac"
$ws.Cells.Item(13, 2).Value = "Gnai2"
$ws.Cells.Item(13, 3).Value = "P2ry12"
$ws.Cells.Item(13, 4).Value = "Resolving-Mac"
$ws.Cells.Item(13, 5).Value = 3
$ws.Cells.Item(13, 6).Value = 1
$ws.Cells.Item(13, 7).Value = 89.46554166666668
$ws.Cells.Item(13, 8).Value = 268.396625
$ws.Cells.Item(13, 9).Value = 0.2146705068098712
$ws.Cells.Item(13, 10).Value = 0.2146705068098712
$ws.Cells.Item(13, 11).Value = 3
$ws.Cells.Item(13, 12).Value = 1
$ws.Cells.Item(13, 13).Value = 22.15150833333333
$ws.Cells.Item(13, 14).Value = 66.45452499999999
$ws.Cells.Item(13, 15).Value = 0.987200381317136
$ws.Cells.Item(13, 16).Value = 0.987200381317136
$ws.Cells.Item(13, 17).Value = 1981.796691775347
$ws.Cells.Item(13, 18).Value = 17836.17022597812
$ws.Cells.Item(13, 19).Value = 0.2119228061802477
$ws.Cells.Item(13, 20).Value = 0.2119228061802477
